$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 3051
$ws.Range("I3").Value = 3116
$ws.Range("I4").Value = 745
$ws.Range("I5").Value = 284
$ws.Range("I6").Value = 3567
$ws.Range("I7").Value = 10763

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I4").Value = 25
$ws.Range("I6").Value = 94
$ws.Range("I7").Value = 344

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 54
$ws.Range("I5").Value = 5
$ws.Range("I7").Value = 196

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 103
$ws.Range("I6").Value = 149
$ws.Range("I7").Value = 431

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 79
$ws.Range("I3").Value = 61
$ws.Range("I7").Value = 238

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I5").Value = 35
$ws.Range("I7").Value = 363
$ws.Range("I8").Value = 675
$ws.Range("I15").Value = 135
$ws.Range("I16").Value = 28
$ws.Range("I18").Value = 75
$ws.Range("I20").Value = 278
$ws.Range("I23").Value = 95
$ws.Range("I29").Value = 711
$ws.Range("I33").Value = 494
$ws.Range("I36").Value = 146
$ws.Range("I37").Value = 344
$ws.Range("I42").Value = 375
$ws.Range("I48").Value = 122
$ws.Range("I50").Value = 50
$ws.Range("I51").Value = 97
$ws.Range("I52").Value = 229
$ws.Range("I53").Value = 121
$ws.Range("I54").Value = 242
$ws.Range("I55").Value = 118
$ws.Range("I57").Value = 36
$ws.Range("I60").Value = 53
$ws.Range("I63").Value = 45
$ws.Range("I64").Value = 96
$ws.Range("I65").Value = 238
$ws.Range("I67").Value = 431
$ws.Range("I74").Value = 25
$ws.Range("I75").Value = 36
$ws.Range("I76").Value = 166
$ws.Range("I79").Value = 273
$ws.Range("I82").Value = 12
$ws.Range("I83").Value = 215
$ws.Range("I85").Value = 490
$ws.Range("I86").Value = 63
$ws.Range("I87").Value = 17
$ws.Range("I88").Value = 97
$ws.Range("I90").Value = 131
$ws.Range("I92").Value = 35
$ws.Range("I95").Value = 170
$ws.Range("I99").Value = 196
$ws.Range("I101").Value = 10763

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 76
$ws.Range("I3").Value = 81
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 215

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 58
$ws.Range("I3").Value = 66
$ws.Range("I7").Value = 170

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 118
$ws.Range("I6").Value = 161
$ws.Range("I7").Value = 494

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 121
$ws.Range("I7").Value = 242

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 220
$ws.Range("I3").Value = 249
$ws.Range("I4").Value = 28
$ws.Range("I6").Value = 189
$ws.Range("I7").Value = 711

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I6").Value = 70
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I3").Value = 41
$ws.Range("I7").Value = 166

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 127
$ws.Range("I3").Value = 198
$ws.Range("I4").Value = 27
$ws.Range("I7").Value = 490

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 102
$ws.Range("I3").Value = 126
$ws.Range("I5").Value = 14
$ws.Range("I7").Value = 375

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 31
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I2").Value = 26
$ws.Range("I7").Value = 95

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 81
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 273

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I3").Value = 31
$ws.Range("I7").Value = 96

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 77
$ws.Range("I4").Value = 18
$ws.Range("I5").Value = 9
$ws.Range("I7").Value = 278

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 45
$ws.Range("I6").Value = 45
$ws.Range("I7").Value = 146

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 83
$ws.Range("I7").Value = 229

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I3").Value = 32
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I2").Value = 13
$ws.Range("I4").Value = 11
$ws.Range("I7").Value = 50

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("I2").Value = 12
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I3").Value = 33
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 210
$ws.Range("I3").Value = 185
$ws.Range("I6").Value = 217
$ws.Range("I7").Value = 675

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I3").Value = 10
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I6").Value = 10
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 40
$ws.Range("I3").Value = 26
$ws.Range("I7").Value = 131

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I2").Value = 18
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I3").Value = 15
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("I5").Value = 7
$ws.Range("I6").Value = 12

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 105
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 363

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("I6").Value = 8
$ws.Range("I7").Value = 17

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("I3").Value = 1
$ws.Range("I7").Value = 25
